$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.173.24"
$ws.Range("E2").Value = "  -1.46%  "

$ws.Range("D3").Value = "2.912.48"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'348.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").Value = "'106.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.12%  "

$ws.Range("D7").Value = "'0.549"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.33%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("E10").Value = "  -5.03%  "

$ws.Range("E11").Value = "  +1.34%  "

$ws.Range("E12").Value = "  -3.52%  "

$ws.Range("D13").Value = "'18.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.83%  "

$ws.Range("D14").Value = "3.366.39"
$ws.Range("E14").Value = "  +0.23%  "

$ws.Range("D15").Value = "'7.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").Value = "2.919.59"
$ws.Range("E16").Value = "  +0.56%  "

$ws.Range("D18").Value = "51.251.22"
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").Value = "'3.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.65%  "

$ws.Range("D20").Value = "'7.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.97%  "

$ws.Range("D21").Value = "'13.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.78%  "

$ws.Range("D22").Value = "0.0₃0957"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").Value = "'68.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.18%  "

$ws.Range("D24").Value = "'260.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.80%  "

$ws.Range("D25").Value = "'2.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.47%  "

$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.170"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.21%  "

$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D27").Value = "'7.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.13%  "

$ws.Range("D28").Value = "'26.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.79%  "

$ws.Range("E29").Value = "  +0.01%  "

$ws.Range("D30").Value = "'0.104"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.45%  "

$ws.Range("D31").Value = "'10.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.55%  "

$ws.Range("D32").Value = "'6.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.23%  "

$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'35.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.33%  "

$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.99%  "

$ws.Range("D35").Value = "'50.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.33%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").Value = "'0.0423"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.26%  "

$ws.Range("D38").Value = "'3.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.30%  "

$ws.Range("D39").Value = "'17.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.45%  "

$ws.Range("E40").Value = "  -4.77%  "

$ws.Range("D41").Value = "'2.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "

$ws.Range("E42").Value = "  -1.00%  "

$ws.Range("D43").Value = "'22.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.75%  "

$ws.Range("D44").Value = "'118.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.40%  "

$ws.Range("D45").Value = "'2.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.26%  "

$ws.Range("D46").Value = "2.082.65"
$ws.Range("E46").Value = "  -4.70%  "

$ws.Range("D47").Value = "'3.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.55%  "

$ws.Range("D48").Value = "'2.22"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.37%  "

$ws.Range("D49").Value = "'0.237"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.38%  "

$ws.Range("E50").Value = "  -1.41%  "

$ws.Range("D51").Value = "'0.887"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.75%  "
